$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D so the existing "Tipo" column (D) shifts to E
$ws.Columns.Item(4).Insert()

# New header "MAE" in D1, matching the style used by the other header cells (copy from C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D1").Value = "MAE"

# Update existing values B2 and C2
$ws.Range("B2").Value = 0.09089304155458464
$ws.Range("C2").Value = 0.9987693643238456

# New numeric value for MAE in D2
$ws.Range("D2").Value = 0.237869009368621

$wb.Save()
